$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) tweaks
$ws.Range("B2").Value = 260.03792218044629
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 255.93217930074152
$ws.Range("E2").Value = 306.61340535368265

# Row 3 (STR) tweaks
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 302.58092235049878
$ws.Range("D3").Value = 249.04615181701169
$ws.Range("E3").Value = 308.11656894531689

# Selection narrowed to the touched range
$ws.Range("B1:E3").Select()
